$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Update the "Vaisselle" challenge text on the "defis" sheet (row 1, column B)
# to add the "pendant 364 jours" emphasis.
$newText = "Vaisselle du quotidien : garde le nombre nécessaire pour tenir entre deux vaisselles<br>Vaisselle pour les occasions : en as-tu réellement besoin ? Un service que tu utilises une fois par an ne sert à rien pendant 364 jours...<br>Cassé/abîmé : à jeter ou recycler."
$ws2.Range("B1").Value = $newText

# The longer text now wraps onto a third line, so the row needs to grow.
$ws2.Rows.Item(1).RowHeight = 43.2

# The "defis" sheet becomes the active / selected sheet and tab, with B7
# as the last selected cell (instead of "parcours" / B9 before).
$ws2.Activate()
$ws2.Range("B7").Select()
